$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New filter fields added alongside the existing "Bookshelves" category value:
#   B1 -> category filter, reuses the same "Bookshelves" text as A1
#   C1 -> PriceRange filter, numeric 15000
#   D1 -> status filter, "Open"
$ws.Range("B1").Value = "Bookshelves"
$ws.Range("C1").Value = 15000
$ws.Range("D1").Value = "Open"

# Set page orientation to portrait (adds <pageSetup .../> to the sheet).
$ws.PageSetup.Orientation = 1

# Move/restore the active selection.
$ws.Range("F9").Select()
